{"js": "const replacements = [\n  [\"909\u00f76=\", \"458\u00f72=\"],\n  [\"220\u00f73=\", \"511\u00f75=\"],\n  [\"568\u00f76=\", \"150\u00f73=\"],\n  [\"276\u00f72=\", \"941\u00f73=\"],\n  [\"234\u00f76=\", \"920\u00f75=\"],\n  [\"288\u00f78=\", \"217\u00f76=\"],\n  [\"882\u00f75=\", \"967\u00f74=\"],\n  [\"525\u00f75=\", \"384\u00f77=\"],\n  [\"121\u00f72=\", \"597\u00f76=\"],\n  [\"647\u00f78=\", \"306\u00f76=\"],\n  [\"308\u00f75=\", \"562\u00f75=\"],\n  [\"565\u00f72=\", \"494\u00f77=\"],\n  [\"780\u00f72=\", \"459\u00f76=\"],\n  [\"926\u00f74=\", \"265\u00f79=\"],\n  [\"515\u00f73=\", \"794\u00f77=\"],\n  [\"935\u00f73=\", \"363\u00f75=\"],\n  [\"809\u00f74=\", \"657\u00f73=\"],\n  [\"278\u00f73=\", \"985\u00f73=\"],\n  [\"326\u00f78=\", \"766\u00f73=\"],\n  [\"751\u00f77=\", \"407\u00f79=\"],\n  [\"468\u00f76=\", \"206\u00f78=\"],\n  [\"961\u00f77=\", \"151\u00f72=\"],\n  [\"768\u00f74=\", \"185\u00f72=\"],\n  [\"915\u00f79=\", \"369\u00f77=\"],\n  [\"233\u00f74=\", \"210\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"909\u00f76=\", \"458\u00f72=\"),\n    @(\"220\u00f73=\", \"511\u00f75=\"),\n    @(\"568\u00f76=\", \"150\u00f73=\"),\n    @(\"276\u00f72=\", \"941\u00f73=\"),\n    @(\"234\u00f76=\", \"920\u00f75=\"),\n    @(\"288\u00f78=\", \"217\u00f76=\"),\n    @(\"882\u00f75=\", \"967\u00f74=\"),\n    @(\"525\u00f75=\", \"384\u00f77=\"),\n    @(\"121\u00f72=\", \"597\u00f76=\"),\n    @(\"647\u00f78=\", \"306\u00f76=\"),\n    @(\"308\u00f75=\", \"562\u00f75=\"),\n    @(\"565\u00f72=\", \"494\u00f77=\"),\n    @(\"780\u00f72=\", \"459\u00f76=\"),\n    @(\"926\u00f74=\", \"265\u00f79=\"),\n    @(\"515\u00f73=\", \"794\u00f77=\"),\n    @(\"935\u00f73=\", \"363\u00f75=\"),\n    @(\"809\u00f74=\", \"657\u00f73=\"),\n    @(\"278\u00f73=\", \"985\u00f73=\"),\n    @(\"326\u00f78=\", \"766\u00f73=\"),\n    @(\"751\u00f77=\", \"407\u00f79=\"),\n    @(\"468\u00f76=\", \"206\u00f78=\"),\n    @(\"961\u00f77=\", \"151\u00f72=\"),\n    @(\"768\u00f74=\", \"185\u00f72=\"),\n    @(\"915\u00f79=\", \"369\u00f77=\"),\n    @(\"233\u00f74=\", \"210\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}"}
